$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 1927
$ws.Range("I70").Value = 1851.5
$ws.Range("J70").Value = 1964.75
$ws.Range("K70").Value = 5554.5
$ws.Range("L70").Value = 5894.25
$ws.Range("M70").Value = -5284.5
$ws.Range("N70").Value = -6434.25
# Row 73
$ws.Range("H73").Value = 1927
$ws.Range("I73").Value = 1851.5
$ws.Range("J73").Value = 1964.75
$ws.Range("K73").Value = 5554.5
$ws.Range("L73").Value = 5894.25
$ws.Range("M73").Value = -4618.5
$ws.Range("N73").Value = -7766.25
# Row 107
$ws.Range("H107").Value = 556068.25
$ws.Range("I107").Value = 741051.1
$ws.Range("J107").Value = 1119.6
$ws.Range("K107").Value = 741051.1
$ws.Range("L107").Value = 1119.6
$ws.Range("M107").Value = -739131.1
$ws.Range("N107").Value = -4959.6
# Row 116
$ws.Range("H116").Value = 4327448.5
$ws.Range("I116").Value = 23063166
$ws.Range("J116").Value = 3821.7693
$ws.Range("K116").Value = 23063166
$ws.Range("L116").Value = 3821.7693
$ws.Range("M116").Value = -23059724
$ws.Range("N116").Value = -10705.7693
# Row 132
$ws.Range("H132").Value = 295452.12
$ws.Range("I132").Value = 419784.7
$ws.Range("J132").Value = 37906.07
$ws.Range("K132").Value = 1259354.1
$ws.Range("L132").Value = 113718.21
$ws.Range("M132").Value = -1256824.1
$ws.Range("N132").Value = -118778.21
# Row 137
$ws.Range("H137").Value = 1687.2609
$ws.Range("I137").Value = 1731.7894
$ws.Range("J137").Value = 1475.75
$ws.Range("K137").Value = 5195.3682
$ws.Range("L137").Value = 4427.25
$ws.Range("M137").Value = -2645.3682
$ws.Range("N137").Value = -9527.25
# Row 138
$ws.Range("H138").Value = 1680.8
$ws.Range("I138").Value = 771.2308
$ws.Range("J138").Value = 2000.3784
$ws.Range("K138").Value = 2313.6924
$ws.Range("L138").Value = 6001.135200000001
$ws.Range("M138").Value = 2826.3076
$ws.Range("N138").Value = -16281.1352

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2072.8472
$ws.Range("I32").Value = 1267.8793
$ws.Range("J32").Value = 5407.7144
$ws.Range("K32").Value = 1267.8793
$ws.Range("L32").Value = 5407.7144
$ws.Range("M32").Value = -980.8793000000001
$ws.Range("N32").Value = -5981.7144
# Row 74
$ws.Range("H74").Value = 6952.6665
$ws.Range("I74").Value = 1337.6666
$ws.Range("K74").Value = 1337.6666
$ws.Range("M74").Value = -463.6666
# Row 77
$ws.Range("H77").Value = 6952.6665
$ws.Range("I77").Value = 1337.6666
$ws.Range("K77").Value = 6688.333000000001
$ws.Range("M77").Value = -2320.333000000001
# Row 102
$ws.Range("H102").Value = 4301.6665
$ws.Range("I102").Value = 4562
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 4562
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -2940
$ws.Range("N102").Value = -6244
# Row 122
$ws.Range("H122").Value = 9714.333000000001
$ws.Range("I122").Value = 13463.75
$ws.Range("J122").Value = 2215.5
$ws.Range("K122").Value = 40391.25
$ws.Range("L122").Value = 6646.5
$ws.Range("M122").Value = -37941.25
$ws.Range("N122").Value = -11546.5
# Row 133
$ws.Range("H133").Value = 64950
$ws.Range("J133").Value = 64950
$ws.Range("L133").Value = 64950
$ws.Range("N133").Value = -70010

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1532.8846
$ws.Range("I20").Value = 1283.8235
$ws.Range("J20").Value = 2003.3334
$ws.Range("K20").Value = 1283.8235
$ws.Range("L20").Value = 2003.3334
$ws.Range("M20").Value = -1036.8235
$ws.Range("N20").Value = -2497.3334
# Row 105
$ws.Range("H105").Value = 2505.7297
$ws.Range("I105").Value = 2748.1904
$ws.Range("K105").Value = 2748.1904
$ws.Range("M105").Value = -1001.1904
# Row 134
$ws.Range("H134").Value = 4136.0527
$ws.Range("I134").Value = 3294.818
$ws.Range("K134").Value = 9884.454000000002
$ws.Range("M134").Value = -7349.454000000002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 580.2
$ws.Range("I22").Value = 150.25
$ws.Range("J22").Value = 2300
$ws.Range("K22").Value = 150.25
$ws.Range("L22").Value = 2300
$ws.Range("M22").Value = 199.75
$ws.Range("N22").Value = -3000
# Row 31
$ws.Range("H31").Value = 4926.95
$ws.Range("I31").Value = 1851.375
$ws.Range("J31").Value = 6977.3335
$ws.Range("K31").Value = 1851.375
$ws.Range("L31").Value = 6977.3335
$ws.Range("M31").Value = -1556.375
$ws.Range("N31").Value = -7567.3335
# Row 34
$ws.Range("H34").Value = 4926.95
$ws.Range("I34").Value = 1851.375
$ws.Range("J34").Value = 6977.3335
$ws.Range("K34").Value = 1851.375
$ws.Range("L34").Value = 6977.3335
$ws.Range("M34").Value = -1649.375
$ws.Range("N34").Value = -7381.3335
# Row 134
$ws.Range("H134").Value = 3517.88
$ws.Range("I134").Value = 2097.7222
$ws.Range("K134").Value = 6293.1666
$ws.Range("M134").Value = -3758.1666

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 687.3333
$ws.Range("I86").Value = 731
$ws.Range("J86").Value = 600
$ws.Range("K86").Value = 2193
$ws.Range("L86").Value = 1800
$ws.Range("M86").Value = -1007
$ws.Range("N86").Value = -4172
# Row 89
$ws.Range("H89").Value = 687.3333
$ws.Range("I89").Value = 731
$ws.Range("J89").Value = 600
$ws.Range("K89").Value = 6579
$ws.Range("L89").Value = 5400
$ws.Range("M89").Value = -651
$ws.Range("N89").Value = -17256
# Row 107
$ws.Range("H107").Value = 293.33334
$ws.Range("I107").Value = 290
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 870
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 1050
$ws.Range("N107").Value = -4740
# Row 129
$ws.Range("H129").Value = 1106.2307
$ws.Range("J129").Value = 1347.8889
$ws.Range("L129").Value = 4043.6667
$ws.Range("N129").Value = -14043.6667
# Row 132
$ws.Range("H132").Value = 969.2632
$ws.Range("I132").Value = 961.6
$ws.Range("K132").Value = 8654.4
$ws.Range("M132").Value = -6124.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2725
$ws.Range("J80").Value = 2737.5
$ws.Range("L80").Value = 2737.5
$ws.Range("N80").Value = -4733.5
# Row 83
$ws.Range("H83").Value = 2725
$ws.Range("J83").Value = 2737.5
$ws.Range("L83").Value = 13687.5
$ws.Range("N83").Value = -23671.5
# Row 107
$ws.Range("H107").Value = 847.6
$ws.Range("J107").Value = 510
$ws.Range("L107").Value = 510
$ws.Range("N107").Value = -4350
# Row 113
$ws.Range("H113").Value = 3500
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -9340
# Row 122
$ws.Range("H122").Value = 1113951.8
$ws.Range("I122").Value = 1854501.6
$ws.Range("K122").Value = 5563504.800000001
$ws.Range("M122").Value = -5561054.800000001
# Row 140
$ws.Range("H140").Value = 51832.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 51832.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 51832.5
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -62192.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 56536.777
$ws.Range("I82").Value = 63353.875
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 63353.875
$ws.Range("L82").Value = 2000
$ws.Range("M82").Value = -62992.875
$ws.Range("N82").Value = -2722
# Row 85
$ws.Range("H85").Value = 56536.777
$ws.Range("I85").Value = 63353.875
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 63353.875
$ws.Range("L85").Value = 2000
$ws.Range("M85").Value = -62105.875
$ws.Range("N85").Value = -4496
# Row 94
$ws.Range("H94").Value = 24999.666
$ws.Range("J94").Value = 24999.666
$ws.Range("L94").Value = 24999.666
$ws.Range("N94").Value = -26351.666
# Row 136
$ws.Range("H136").Value = 2898.9119
$ws.Range("I136").Value = 1604.9231
$ws.Range("J136").Value = 3699.9524
$ws.Range("K136").Value = 4814.7693
$ws.Range("L136").Value = 11099.8572
$ws.Range("M136").Value = -2264.7693
$ws.Range("N136").Value = -16199.8572

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 57278.95
$ws.Range("I81").Value = 1000000
$ws.Range("J81").Value = 4905.5557
$ws.Range("K81").Value = 2000000
$ws.Range("L81").Value = 9811.1114
$ws.Range("M81").Value = -1998939
$ws.Range("N81").Value = -11933.1114
# Row 84
$ws.Range("H84").Value = 57278.95
$ws.Range("I84").Value = 1000000
$ws.Range("J84").Value = 4905.5557
$ws.Range("K84").Value = 10000000
$ws.Range("L84").Value = 49055.557
$ws.Range("M84").Value = -9994696
$ws.Range("N84").Value = -59663.557
